# feat: add 2022-Q4 data
#
# 1. Insert a new "2022-Q4" sheet right after "总计", populated with its
#    fund-holding data (copied/grown from the "2022-Q1" sheet template so it
#    inherits the exact same styles/number formats).
# 2. Update the "总计" (totals) sheet: shift the existing quarter rows down
#    by one and insert a new row for 2022-Q4 at the top of the data, then
#    renumber the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: build the new "2022-Q4" worksheet
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)           # "2022-Q1" sheet (same header set)
$template.Copy($null, $wb.Worksheets.Item(1))
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# the template only carries 2 data rows (rows 2-3); grow it to 6 data rows
# (rows 2-7) by duplicating row 3's formatting downward.
$q4.Range("A3:H3").Copy($q4.Range("A4:H4"))
$q4.Range("A3:H3").Copy($q4.Range("A5:H5"))
$q4.Range("A3:H3").Copy($q4.Range("A6:H6"))
$q4.Range("A3:H3").Copy($q4.Range("A7:H7"))

# columns B, D, E, F, G hold numeric-looking text (fund codes / percentages)
# that must stay text (leading zeros, fixed decimals) instead of being
# coerced to numbers.
$q4.Range("B2:B7").NumberFormat = "@"
$q4.Range("D2:G7").NumberFormat = "@"

$q4rows = @(
    @("159758", "华夏中证红利质量ETF",              "1.81", "99.33", "3.36", "0.0608", 9),
    @("016174", "汇丰晋信策略优选混合A",              "1.78", "74.92", "2.30", "0.0409", 9),
    @("159628", "万家国证2000ETF",                   "2.55", "97.28", "0.44", "0.0112", 9),
    @("016175", "汇丰晋信策略优选混合C",              "0.40", "74.92", "2.30", "0.0092", 9),
    @("165524", "信诚中证智能家居指数（LOF）A",        "0.37", "91.28", "0.99", "0.0037", 10),
    @("013084", "信诚中证智能家居指数（LOF）C",        "0.14", "91.28", "0.99", "0.0014", 10)
)

for ($i = 0; $i -lt $q4rows.Length; $i++) {
    $r = $i + 2
    $row = $q4rows[$i]
    $q4.Range("A$r").Value = $i
    $q4.Range("B$r").Value = $row[0]
    $q4.Range("C$r").Value = $row[1]
    $q4.Range("D$r").Value = $row[2]
    $q4.Range("E$r").Value = $row[3]
    $q4.Range("F$r").Value = $row[4]
    $q4.Range("G$r").Value = $row[5]
    $q4.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 2: update the "总计" summary sheet with the new 2022-Q4 row
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# shift existing data rows 2-6 down to 3-7 (bottom-up so nothing is clobbered),
# carrying the row's formatting along with it.
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $total.Range("B$r`:D$r").Copy($total.Range("B$dest`:D$dest"))
}

# new top row: 2022-Q4 summary
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.13

# renumber the leading index column for all 6 data rows (0..5)
for ($i = 0; $i -le 5; $i++) {
    $r = $i + 2
    $total.Range("A$r").Value = $i
}

# keep "总计" as the active tab, matching the workbook's original bookViews
# (activeTab="0"), which the edit does not touch.
$total.Activate()
